# Scheduled data refresh: update Leve price/profit columns (H-N) across all job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 749
$ws.Range("I18").Value = 749
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 749
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -465
$ws.Range("N18").ClearContents()

$ws.Range("H51").Value = 23125
$ws.Range("I51").Value = 100000
$ws.Range("J51").Value = 12142.857
$ws.Range("K51").Value = 100000
$ws.Range("L51").Value = 12142.857
$ws.Range("M51").Value = -99516
$ws.Range("N51").Value = -13110.857

$ws.Range("H127").Value = 4032.8333
$ws.Range("I127").Value = 4032.8333
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 12098.4999
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = -7138.499899999999
$ws.Range("N127").ClearContents()

$ws.Range("H137").Value = 7010.0454
$ws.Range("I137").Value = 5984.3335
$ws.Range("K137").Value = 17953.0005
$ws.Range("M137").Value = -15403.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 543.625
$ws.Range("I32").Value = 543.625
$ws.Range("K32").Value = 543.625
$ws.Range("M32").Value = -256.625

$ws.Range("H63").Value = 3478.7778
$ws.Range("I63").Value = 3478.7778
$ws.Range("K63").Value = 3478.7778
$ws.Range("M63").Value = -2792.7778

$ws.Range("H66").Value = 3478.7778
$ws.Range("I66").Value = 3478.7778
$ws.Range("K66").Value = 17393.889
$ws.Range("M66").Value = -13961.889

$ws.Range("H102").Value = 2000
$ws.Range("J102").Value = 2000
$ws.Range("L102").Value = 2000
$ws.Range("N102").Value = -5244

$ws.Range("H122").Value = 2377.75
$ws.Range("I122").Value = 2170.3333
$ws.Range("K122").Value = 6510.999899999999
$ws.Range("M122").Value = -4060.999899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3975.6667
$ws.Range("I20").Value = 3975.6667
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 3975.6667
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -3728.6667
$ws.Range("N20").ClearContents()

$ws.Range("H80").Value = 1751.4
$ws.Range("I80").Value = 1751.4
$ws.Range("K80").Value = 1751.4
$ws.Range("M80").Value = -753.4000000000001

$ws.Range("H83").Value = 1751.4
$ws.Range("I83").Value = 1751.4
$ws.Range("K83").Value = 8757
$ws.Range("M83").Value = -3765

$ws.Range("H94").Value = 967.25
$ws.Range("I94").Value = 853
$ws.Range("K94").Value = 853
$ws.Range("M94").Value = -402

$ws.Range("H107").Value = 1500
$ws.Range("J107").Value = 2000
$ws.Range("L107").Value = 2000
$ws.Range("N107").Value = -5840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7340.3335
$ws.Range("I31").Value = 5181.909
$ws.Range("K31").Value = 5181.909
$ws.Range("M31").Value = -4886.909

$ws.Range("H34").Value = 7340.3335
$ws.Range("I34").Value = 5181.909
$ws.Range("K34").Value = 5181.909
$ws.Range("M34").Value = -4979.909

$ws.Range("H132").Value = 3014.7778
$ws.Range("I132").Value = 1590.6428
$ws.Range("K132").Value = 4771.928400000001
$ws.Range("M132").Value = -2241.928400000001

$ws.Range("H134").Value = 6345.9
$ws.Range("I134").Value = 3682.375
$ws.Range("K134").Value = 11047.125
$ws.Range("M134").Value = -8512.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 89.166664
$ws.Range("I40").Value = 47
$ws.Range("J40").Value = 300
$ws.Range("K40").Value = 188
$ws.Range("L40").Value = 1200
$ws.Range("M40").Value = -119
$ws.Range("N40").Value = -1338

$ws.Range("H52").Value = 2000
$ws.Range("J52").Value = 2000
$ws.Range("L52").Value = 6000
$ws.Range("N52").Value = -6532

$ws.Range("H116").Value = 1070
$ws.Range("I116").Value = 1070
$ws.Range("K116").Value = 3210
$ws.Range("M116").Value = 232

$ws.Range("H134").Value = 60000
$ws.Range("I134").Value = 60000
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 180000
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -174930
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 10000600
$ws.Range("I10").Value = 900
$ws.Range("K10").Value = 900
$ws.Range("M10").Value = -731

$ws.Range("H45").Value = 40000
$ws.Range("J45").Value = 40000
$ws.Range("L45").Value = 40000
$ws.Range("N45").Value = -41118

$ws.Range("H126").Value = 10622
$ws.Range("I126").Value = 10622
$ws.Range("K126").Value = 31866
$ws.Range("M126").Value = -29396

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 10832.167
$ws.Range("J22").Value = 9995.333000000001
$ws.Range("L22").Value = 9995.333000000001
$ws.Range("N22").Value = -10585.333

$ws.Range("H27").Value = 10832.167
$ws.Range("J27").Value = 9995.333000000001
$ws.Range("L27").Value = 9995.333000000001
$ws.Range("N27").Value = -10209.333

$ws.Range("H46").Value = 5561.2583
$ws.Range("I46").Value = 4974.75
$ws.Range("J46").Value = 5648.148
$ws.Range("K46").Value = 4974.75
$ws.Range("L46").Value = 5648.148
$ws.Range("M46").Value = -4786.75
$ws.Range("N46").Value = -6024.148

$ws.Range("H55").Value = 2192.8
$ws.Range("I55").Value = 1833.3334
$ws.Range("J55").Value = 2732
$ws.Range("K55").Value = 1833.3334
$ws.Range("L55").Value = 2732
$ws.Range("M55").Value = -1660.3334
$ws.Range("N55").Value = -3078

$ws.Range("H74").Value = 30000
$ws.Range("J74").Value = 30000
$ws.Range("L74").Value = 30000
$ws.Range("N74").Value = -31996

$ws.Range("H77").Value = 30000
$ws.Range("J77").Value = 30000
$ws.Range("L77").Value = 90000
$ws.Range("N77").Value = -99984

$ws.Range("H93").Value = 1733
$ws.Range("I93").Value = 1879.6
$ws.Range("K93").Value = 1879.6
$ws.Range("M93").Value = -631.5999999999999

$ws.Range("H136").Value = 8706.879999999999
$ws.Range("I136").Value = 3644.5386
$ws.Range("K136").Value = 10933.6158
$ws.Range("M136").Value = -8383.6158

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 500
$ws.Range("I7").Value = 500
$ws.Range("K7").Value = 500
$ws.Range("M7").Value = -387

$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()

$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()

$ws.Range("H81").Value = 8000
$ws.Range("I81").Value = 2000
$ws.Range("K81").Value = 4000
$ws.Range("M81").Value = -2939

$ws.Range("H84").Value = 8000
$ws.Range("I84").Value = 2000
$ws.Range("K84").Value = 20000
$ws.Range("M84").Value = -14696

$ws.Range("H122").Value = 12003.5
$ws.Range("I122").Value = 12003.5
$ws.Range("K122").Value = 36010.5
$ws.Range("M122").Value = -33560.5

$ws.Range("H126").Value = 1800
$ws.Range("I126").Value = 1800
$ws.Range("K126").Value = 5400
$ws.Range("M126").Value = -2930
